$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update VALIDATIONS (column J) for the type-ahead dependency test rows,
# switching the validation field name from "value" to "suggestions.keyword".
$ws.Range("J4").Value = "status=200||suggestions.keyword=mic"
$ws.Range("J5").Value = "status=200||suggestions.keyword=bio||suggestions.keyword=methanol"
$ws.Range("J6").Value = "status=200||suggestions.keyword=bio||suggestions.keyword=methanol"

# Widen column J (VALIDATIONS) to fit the new, longer values and drop the
# "best fit" auto-width flag in favor of an explicit custom width
# (~50.14 characters wide).
$ws.Columns.Item(10).ColumnWidth = 49.3

# Update the sheet view: scroll so column H is the left-most visible column,
# and move the active selection cell to the top of the highlighted range
# (selecting L2:L6 starting from L2 makes L2 the active cell while keeping
# the same highlighted range).
$excel.ActiveWindow.ScrollColumn = 8
$ws.Range("L2:L6").Select()
